$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number-format on cells whose new value would otherwise be
# auto-converted to a numeric literal by Excel (these columns are text-typed
# in the source data, e.g. "581.00", "0.999", "0.0840").
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range('D2').Value = '66.984.81'
$ws.Range('E2').Value = '  -5.50%  '
$ws.Range('D3').Value = '3.223.26'
$ws.Range('E3').Value = '  -8.71%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '581.00'
$ws.Range('E5').Value = '  -5.30%  '
$ws.Range('D6').Value = '151.08'
$ws.Range('E6').Value = '  -13.10%  '
$ws.Range('D8').Value = '3.214.35'
$ws.Range('E8').Value = '  -8.87%  '
$ws.Range('D9').Value = '0.544'
$ws.Range('E9').Value = '  -11.05%  '
$ws.Range('E10').Value = '  -12.68%  '
$ws.Range('D11').Value = '6.83'
$ws.Range('E11').Value = '  -5.39%  '
$ws.Range('D12').Value = '0.503'
$ws.Range('E12').Value = '  -14.68%  '
$ws.Range('D13').Value = '38.13'
$ws.Range('E13').Value = '  -18.18%  '
$ws.Range('E14').Value = '  -12.11%  '
$ws.Range('D15').Value = '3.737.81'
$ws.Range('E15').Value = '  -8.81%  '
$ws.Range('D16').Value = '66.839.60'
$ws.Range('E16').Value = '  -5.69%  '
$ws.Range('D17').Value = '3.221.87'
$ws.Range('E17').Value = '  -9.02%  '
$ws.Range('D18').Value = '541.13'
$ws.Range('E18').Value = '  -11.92%  '
$ws.Range('E19').Value = '  -5.83%  '
$ws.Range('E20').Value = '  -15.50%  '
$ws.Range('D21').Value = '15.07'
$ws.Range('E21').Value = '  -15.15%  '
$ws.Range('E22').Value = '  -14.69%  '
$ws.Range('D23').Value = '7.68'
$ws.Range('E23').Value = '  -14.65%  '
$ws.Range('D24').Value = '85.57'
$ws.Range('E24').Value = '  -12.80%  '
$ws.Range('D25').Value = '13.39'
$ws.Range('E25').Value = '  -14.95%  '
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').Value = '3.14'
$ws.Range('E27').Value = '  -16.84%  '
$ws.Range('D28').Value = '8.06'
$ws.Range('E28').Value = '  -11.99%  '
$ws.Range('D29').Value = '29.33'
$ws.Range('E29').Value = '  -13.45%  '
$ws.Range('D30').Value = '2.12'
$ws.Range('E30').Value = '  -18.37%  '
$ws.Range('E31').Value = '  -15.15%  '
$ws.Range('D32').Value = '1.13'
$ws.Range('E32').Value = '  -12.93%  '
$ws.Range('D33').Value = '549.35'
$ws.Range('E33').Value = '  -9.23%  '
$ws.Range('D34').Value = '6.53'
$ws.Range('E34').Value = '  -20.11%  '
$ws.Range('E35').Value = '  -17.16%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('D37').Value = '53.22'
$ws.Range('E37').Value = '  -6.64%  '
$ws.Range('E38').Value = '  -8.71%  '
$ws.Range('B39').Value = 'Cosmos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D39').Value = '9.14'
$ws.Range('E39').Value = '  -15.82%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = '0.0840'
$ws.Range('E40').Value = '  -16.54%  '
$ws.Range('D41').Value = '0.125'
$ws.Range('E41').Value = '  -14.03%  '
$ws.Range('D42').Value = '2.927.44'
$ws.Range('E42').Value = '  -13.06%  '
$ws.Range('E43').Value = '  -27.69%  '
$ws.Range('E44').Value = '  -16.93%  '
$ws.Range('D45').Value = '0.0₃0584'
$ws.Range('E45').Value = '  -21.43%  '
$ws.Range('E46').Value = '  -19.96%  '
$ws.Range('D48').Value = '25.86'
$ws.Range('E48').Value = '  -19.85%  '
$ws.Range('E49').Value = '  -18.79%  '
$ws.Range('E50').Value = '  -13.26%  '
$ws.Range('D51').Value = '122.99'
$ws.Range('E51').Value = '  -8.21%  '
